$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 275
$ws.Range("D9").Value = 2031
$ws.Range("D14").Value = 2008
$ws.Range("D22").Value = 1695
$ws.Range("D29").Value = 2889
$ws.Range("D30").Value = -1063
$ws.Range("D33").Value = -109
$ws.Range("D40").Value = -66
$ws.Range("D41").Value = 2345
$ws.Range("D42").Value = 643
$ws.Range("D43").Value = -2160
$ws.Range("D45").Value = 2272
$ws.Range("D46").Value = 59
$ws.Range("D48").Value = -20
$ws.Range("D51").Value = -1320
$ws.Range("D52").Value = -32
$ws.Range("D53").Value = 12
$ws.Range("D56").Value = 1142
$ws.Range("D63").Value = -183
$ws.Range("D75").Value = -479
$ws.Range("D76").Value = -6
$ws.Range("D82").Value = 124
$ws.Range("D90").Value = 272
$ws.Range("D91").Value = 415
$ws.Range("D92").Value = 673
$ws.Range("D93").Value = 307
$ws.Range("D97").Value = 408
$ws.Range("D99").Value = 95
$ws.Range("D118").Value = -182
$ws.Range("D124").Value = 170
$ws.Range("D125").Value = -76
$ws.Range("D127").Value = 357
$ws.Range("D155").Value = -2358
$ws.Range("D160").Value = 380
$ws.Range("D181").Value = -54
$ws.Range("D182").Value = -132
$ws.Range("D207").Value = -6629
$ws.Range("D225").Value = -50
$ws.Range("D230").Value = -1650
$ws.Range("D231").Value = -26
$ws.Range("D238").Value = -6600
$ws.Range("D242").Value = -240